$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.230.28'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '2.432.68'
$ws.Range("E3").Value = '  -2.19%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '490.70'
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.60'
$ws.Range("E6").Value = '  +3.23%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.616'
$ws.Range("E7").Value = '  +20.65%  '
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.997'
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = '2.453.05'
$ws.Range("E9").Value = '  -1.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.100'
$ws.Range("E10").Value = '  +1.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.62'
$ws.Range("E11").Value = '  -3.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.337'
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("E13").Value = '  +0.86%  '
$ws.Range("D14").Value = '2.869.60'
$ws.Range("E14").Value = '  -1.85%  '
$ws.Range("D15").Value = '57.236.74'
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.79'
$ws.Range("E16").Value = '  -1.82%  '
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("D18").Value = '2.447.53'
$ws.Range("E18").Value = '  -2.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.78'
$ws.Range("E19").Value = '  +4.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '328.96'
$ws.Range("E20").Value = '  +2.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.99'
$ws.Range("E21").Value = '  -2.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.94'
$ws.Range("E23").Value = '  +0.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '57.90'
$ws.Range("E24").Value = '  -0.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.411'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  -0.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.161'
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("D28").Value = '2.544.60'
$ws.Range("E28").Value = '  -2.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.42'
$ws.Range("E29").Value = '  -3.11%  '
$ws.Range("D30").Value = '0.0₃0788'
$ws.Range("E30").Value = '  -0.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.74'
$ws.Range("E32").Value = '  +2.17%  '
$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.70'
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.53'
$ws.Range("E34").Value = '  +1.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.31'
$ws.Range("E35").Value = '  +1.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.15'
$ws.Range("E36").Value = '  -1.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.72'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.858'
$ws.Range("E38").Value = '  -2.14%  '
$ws.Range("E39").Value = '  +10.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '34.22'
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("E41").Value = '  -0.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.52'
$ws.Range("E42").Value = '  +0.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.996'
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.599'
$ws.Range("E44").Value = '  -2.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0538'
$ws.Range("E45").Value = '  -3.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '267.72'
$ws.Range("E46").Value = '  +0.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.71'
$ws.Range("E47").Value = '  -3.08%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0229'
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.21'
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.64'
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.74'
$ws.Range("E51").Value = '  +14.82%  '

Write-Output "Applied 108 cell updates."